$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pin")

# Row 39: SLOT0_PRSENT2_0 -> POWER_KEY
$ws.Range("C39").Value = "POWER_KEY"

# Row 40: SLOT0_PRSENT2_1 -> RESET_KEY
$ws.Range("C40").Value = "RESET_KEY"

# Row 43: SLOT1_PRSENT2_0 -> LED1, Mode INPUT -> OUTPUT
$ws.Range("C43").Value = "LED1"
$ws.Range("D43").Value = "OUTPUT"

# Row 44: SLOT1_PRSENT2_1 -> LED2, Mode INPUT -> OUTPUT
$ws.Range("C44").Value = "LED2"
$ws.Range("D44").Value = "OUTPUT"

# Row 76: Mode INPUT -> OUTPUT (MCU_ATX_ON)
$ws.Range("D76").Value = "OUTPUT"

# Row 85: NC -> LINK_GPIO22, Mode INPUT -> OUTPUT
$ws.Range("C85").Value = "LINK_GPIO22"
$ws.Range("D85").Value = "OUTPUT"

# Row 86: NC -> LINK_GPIO23, Mode INPUT -> OUTPUT
$ws.Range("C86").Value = "LINK_GPIO23"
$ws.Range("D86").Value = "OUTPUT"

# Row 87: NC -> LINK_GPIO24, Mode INPUT -> OUTPUT
$ws.Range("C87").Value = "LINK_GPIO24"
$ws.Range("D87").Value = "OUTPUT"
